$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11723.6889453199
$ws.Range("F2").Value = 79.2392382746251

$ws.Range("C3").Value = 7746.78304517756
$ws.Range("F3").Value = 153.353456893459

$ws.Range("C4").Value = 7614.00686235375
$ws.Range("F4").Value = 141.840008567942

$ws.Range("C5").Value = 11740.1786447282
$ws.Range("F5").Value = 324.317782450555

$ws.Range("C6").Value = 11853.6623139506
$ws.Range("F6").Value = 322.460126360758

$ws.Range("C7").Value = 11527.8329886757
$ws.Range("F7").Value = 323.597174832601
